$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1500
$ws.Range("I20").Value = 1500
$ws.Range("K20").Value = 1500
$ws.Range("M20").Value = -1270
$ws.Range("H35").Value = 1500
$ws.Range("I35").Value = 1500
$ws.Range("K35").Value = 1500
$ws.Range("M35").Value = -1121
$ws.Range("H103").Value = 45455224
$ws.Range("I103").Value = 897
$ws.Range("J103").Value = 71429130
$ws.Range("K103").Value = 2691
$ws.Range("L103").Value = 214287390
$ws.Range("M103").Value = -2105
$ws.Range("N103").Value = -214288562
$ws.Range("H137").Value = 1239.2354
$ws.Range("I137").Value = 955.4167
$ws.Range("J137").Value = 1920.4
$ws.Range("K137").Value = 2866.2501
$ws.Range("L137").Value = 5761.200000000001
$ws.Range("M137").Value = -316.2501000000002
$ws.Range("N137").Value = -10861.2
$ws.Range("H138").Value = 3337.6326
$ws.Range("I138").Value = 2416.8462
$ws.Range("J138").Value = 4378.522
$ws.Range("K138").Value = 7250.5386
$ws.Range("L138").Value = 13135.566
$ws.Range("M138").Value = -2110.5386
$ws.Range("N138").Value = -23415.566

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1265.75
$ws.Range("I74").Value = 779.875
$ws.Range("J74").Value = 2237.5
$ws.Range("K74").Value = 779.875
$ws.Range("L74").Value = 2237.5
$ws.Range("M74").Value = 94.125
$ws.Range("N74").Value = -3985.5
$ws.Range("H77").Value = 1265.75
$ws.Range("I77").Value = 779.875
$ws.Range("J77").Value = 2237.5
$ws.Range("K77").Value = 3899.375
$ws.Range("L77").Value = 11187.5
$ws.Range("M77").Value = 468.625
$ws.Range("N77").Value = -19923.5
$ws.Range("H109").Value = 33331.418
$ws.Range("J109").Value = 33331.418
$ws.Range("L109").Value = 33331.418
$ws.Range("N109").Value = -36105.418

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 1495.4
$ws.Range("I39").Value = 939.3333
$ws.Range("J39").Value = 6500
$ws.Range("K39").Value = 939.3333
$ws.Range("L39").Value = 6500
$ws.Range("M39").Value = -548.3333
$ws.Range("N39").Value = -7282
$ws.Range("H49").Value = 1495.4
$ws.Range("I49").Value = 939.3333
$ws.Range("J49").Value = 6500
$ws.Range("K49").Value = 939.3333
$ws.Range("L49").Value = 6500
$ws.Range("M49").Value = -757.3333
$ws.Range("N49").Value = -6864
$ws.Range("H132").Value = 1628.7812
$ws.Range("I132").Value = 1200.3214
$ws.Range("J132").Value = 4628
$ws.Range("K132").Value = 3600.9642
$ws.Range("L132").Value = 13884
$ws.Range("M132").Value = -1070.9642
$ws.Range("N132").Value = -18944

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 957.6667
$ws.Range("I5").Value = 503.54544
$ws.Range("J5").Value = 1671.2858
$ws.Range("K5").Value = 1510.63632
$ws.Range("L5").Value = 5013.857400000001
$ws.Range("M5").Value = -1398.63632
$ws.Range("N5").Value = -5237.857400000001
$ws.Range("H107").Value = 426.36365
$ws.Range("J107").Value = 478.75
$ws.Range("L107").Value = 1436.25
$ws.Range("N107").Value = -5276.25
$ws.Range("H114").Value = 532818.1
$ws.Range("I114").Value = 10242.272
$ws.Range("J114").Value = 1251359.9
$ws.Range("K114").Value = 30726.816
$ws.Range("L114").Value = 3754079.7
$ws.Range("M114").Value = -27472.816
$ws.Range("N114").Value = -3760587.7
$ws.Range("H131").Value = 794.47
$ws.Range("I131").Value = 407.6154
$ws.Range("J131").Value = 852.2759
$ws.Range("K131").Value = 1222.8462
$ws.Range("L131").Value = 2556.8277
$ws.Range("M131").Value = 3817.1538
$ws.Range("N131").Value = -12636.8277
$ws.Range("H135").Value = 957.6667
$ws.Range("I135").Value = 503.54544
$ws.Range("J135").Value = 1671.2858
$ws.Range("K135").Value = 4531.90896
$ws.Range("L135").Value = 15041.5722
$ws.Range("M135").Value = -1996.90896
$ws.Range("N135").Value = -20111.5722

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 22998.572
$ws.Range("I57").Value = 10000
$ws.Range("J57").Value = 25165
$ws.Range("K57").Value = 10000
$ws.Range("L57").Value = 25165
$ws.Range("M57").Value = -9180
$ws.Range("N57").Value = -26805
$ws.Range("H122").Value = 2268.1538
$ws.Range("I122").Value = 2543.6
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 7630.799999999999
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = -5180.799999999999
$ws.Range("N122").Value = -8950

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2046.742
$ws.Range("I68").Value = 1981.25
$ws.Range("J68").Value = 2116.6
$ws.Range("K68").Value = 1981.25
$ws.Range("L68").Value = 2116.6
$ws.Range("M68").Value = -1232.25
$ws.Range("N68").Value = -3614.6
$ws.Range("H71").Value = 2046.742
$ws.Range("I71").Value = 1981.25
$ws.Range("J71").Value = 2116.6
$ws.Range("K71").Value = 9906.25
$ws.Range("L71").Value = 10583
$ws.Range("M71").Value = -6162.25
$ws.Range("N71").Value = -18071

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6136.875
$ws.Range("I62").Value = 6482.1816
$ws.Range("J62").Value = 5377.2
$ws.Range("K62").Value = 6482.1816
$ws.Range("L62").Value = 5377.2
$ws.Range("M62").Value = -5858.1816
$ws.Range("N62").Value = -6625.2
$ws.Range("H65").Value = 6136.875
$ws.Range("I65").Value = 6482.1816
$ws.Range("J65").Value = 5377.2
$ws.Range("K65").Value = 32410.908
$ws.Range("L65").Value = 26886
$ws.Range("M65").Value = -29290.908
$ws.Range("N65").Value = -33126
$ws.Range("H81").Value = 2211.6296
$ws.Range("I81").Value = 1224.4762
$ws.Range("J81").Value = 5666.6665
$ws.Range("K81").Value = 2448.9524
$ws.Range("L81").Value = 11333.333
$ws.Range("M81").Value = -1387.9524
$ws.Range("N81").Value = -13455.333
$ws.Range("H84").Value = 2211.6296
$ws.Range("I84").Value = 1224.4762
$ws.Range("J84").Value = 5666.6665
$ws.Range("K84").Value = 12244.762
$ws.Range("L84").Value = 56666.665
$ws.Range("M84").Value = -6940.762000000001
$ws.Range("N84").Value = -67274.66500000001
$ws.Range("H122").Value = 1415.7778
$ws.Range("I122").Value = 1305.2667
$ws.Range("K122").Value = 3915.800099999999
$ws.Range("M122").Value = -1465.800099999999
$ws.Range("H132").Value = 15387395
$ws.Range("I132").Value = 18868960
$ws.Range("J132").Value = 10485.25
$ws.Range("K132").Value = 56606880
$ws.Range("L132").Value = 31455.75
$ws.Range("M132").Value = -56604350
$ws.Range("N132").Value = -36515.75

Write-Host "Applied all Durandal Profits updates"